$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.948.35'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '3.687.97'
$ws.Range('E3').Value = '  +0.42%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '237.01'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -2.23%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.90'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.58%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '657.45'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -0.44%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.423'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '1.06'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('D11').Value = '3.686.29'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('B12').Value = 'ShibaInu'
$ws.Range('C12').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0000307'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +14.29%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.209'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +1.92%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '44.01'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -3.09%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '6.76'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').Value = '4.376.29'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '96.770.77'
$ws.Range('E17').Value = '  +0.48%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '9.16'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +3.29%  '
$ws.Range('D19').Value = '3.694.74'
$ws.Range('E19').Value = '  +0.69%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '13.07'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +0.71%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '18.69'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +1.96%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.504'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -4.02%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '520.89'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.61%  '
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('E25').Value = '  +3.49%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '6.94'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '101.17'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('B28').Value = 'Hedera'
$ws.Range('C28').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '0.195'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +17.24%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '13.49'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.69%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '12.50'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.18%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '3.04'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -0.06%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +0.04%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.189'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('E34').Value = '  +2.37%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '656.03'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '32.08'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -2.79%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.592'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.44%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '8.85'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.162'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.01%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '6.81'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +7.57%  '
$ws.Range('E43').Value = '  +4.03%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '40.45'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -11.14%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.486'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +11.20%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.965'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.10%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0464'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.85%  '
$ws.Range('E48').Value = '  +0.45%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '23.61'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -0.10%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '8.74'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.73%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.50'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.66%  '
